$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws1.Range("C2").Value = "approach"

$ws2 = $wb.Worksheets.Item("computational_comparison")
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"
$ws2.Range("G5").ClearContents()
